# Applies the two content edits on slide 6 ("Preprocessing: Grouping registers"):
#   1. Table 6 header cell "Waste" -> "Difference"
#   2. Text Box 8's second line "on the date" -> "on the same date", plus the
#      box is repositioned/resized to fit the longer text (left shifts left,
#      width grows; top/height unchanged).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- 1. Table header: "Waste" -> "Difference" -------------------------------
$tbl = $s.Shapes.Item(3).Table
$cell = $tbl.Cell(1, 7)
$cellRun = $cell.Shape.TextFrame.TextRange.Paragraphs(1).Runs(1)
$cellRun.Text = "Difference"

# --- 2. Text Box 8: update wording and resize/reposition -------------------
$tb = $s.Shapes.Item(5)

$tbRun = $tb.TextFrame.TextRange.Paragraphs(2).Runs(1)
$tbRun.Text = "on the same date"

# Reposition/resize (values are precise point equivalents of the target EMUs
# 3868103 / 1652905 / 491490, chosen so the COM layer's Single(32-bit float)
# truncation reproduces the exact EMU values).
$tb.Left = 304.5750427246094
$tb.Width = 130.15000915527344
$tb.Height = 38.70000076293945
